$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.1417585925264965
$ws.Range("E2").Value = 10.56343377223308
$ws.Range("F2").Value = 32.37482745422905
